# "Generate Report for handoff"
#
# For both localized-language sheets ("zh-cn" and "de-de") the handoff
# report now reflects a failed/ignored handoff instead of a completed one:
#   - Status (B2): "Ready for handoff" -> "Handoff transform failed"
#   - Latest Handoff File (C2): cleared - no handoff file was produced,
#     so the hyperlink + cell are removed entirely
#   - Latest Handoff Datetime (D2): reset to the "never happened" sentinel
#     date "0001-01-01 00:00:00"
#   - Handoff Reason (H2): "Include" -> "Ignored"
#
# The "Overview" sheet shares the same "Ready for handoff" string in its
# own Status column (B2), so its label flips to the new text too.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Overview").Range("B2").Value = "Handoff transform failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find (but don't delete yet - mutating the Hyperlinks collection while
    # iterating it is unsafe) the hyperlink anchored on C2, the
    # "Latest Handoff File" link to the generated .xlf handoff file.
    $c2Link = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $c2Link = $hl
        }
    }
    if ($c2Link -ne $null) {
        $c2Link.Delete()
    }

    # Drop the now-stale handoff file reference entirely (value + format),
    # matching the cell being removed from the sheet.
    $ws.Range("C2").Clear()

    $ws.Range("B2").Value = "Handoff transform failed"
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
